# Commit: "remove the other volatiles (from stem trichomes)"
#
# The "Other volatiles" table (rows 59-90, columns A-F) on the
# "volatiles_candidates" sheet contained retention-time / Kovats-index
# data (and "NA" placeholders) for volatiles coming from stem
# trichomes. This data is removed here, while the cell formatting
# (styles/borders) of that block is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("volatiles_candidates")

# Clear the values/formulas of the "Other volatiles" data block,
# keeping the existing cell styles intact.
$dataRange = $ws.Range("A59:F90")
$dataRange.ClearContents()

# Reflect the new scroll position / selection left after removing
# that block.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 47
$excel.ActiveWindow.ScrollColumn = 1
$dataRange.Select()
